$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new header cell A1 = "Category", matching the style of the rest
# of the header row (B1:W1 use style index 1 - bold/bordered/centered).
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The category column body (A2:A46) previously carried the header style;
# drop that formatting so those cells use the default style.
$ws.Range("A2:A46").ClearFormats()
